$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the NUFUS (population) column for USA, France, England and Turkey
# with their real population figures. The cells are text-typed (shared
# string) values in the original workbook, so force a text number format
# before writing the digit-only strings, then restore the default "Normal"
# style so formatting stays identical to the source cells.
$ws.Range("C2:C5").NumberFormat = "@"

$ws.Range("C2").Value = "712816"
$ws.Range("C3").Value = "2161000"
$ws.Range("C4").Value = "8982000"
$ws.Range("C5").Value = "5663000"

$ws.Range("C2:C5").Style = "Normal"
